# Correct municipales election dates in column B.
# Two serial dates were recorded incorrectly and need to be corrected:
#   39591 (2008-05-23) -> 39516 (2008-03-09)
#   41723 (2014-03-25) -> 41721 (2014-03-23)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2

    if ($val -eq 39591) {
        $cell.Value = 39516
    }
    elseif ($val -eq 41723) {
        $cell.Value = 41721
    }
}
